$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 607
$ws1.Range("F4").Value = 510
$ws1.Range("F5").Value = 506
$ws1.Range("F6").Value = 289
$ws1.Range("F7").Value = 2588
$ws1.Range("F8").Value = 442
$ws1.Range("F9").Value = 7045
$ws1.Range("F10").Value = 188
$ws1.Range("F11").Value = 444
$ws1.Range("F12").Value = 9
$ws1.Range("F13").Value = 102

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 607
$ws4.Range("F4").Value = 510
$ws4.Range("F5").Value = 506
$ws4.Range("F6").Value = 289
$ws4.Range("F9").Value = 2588
$ws4.Range("F10").Value = 442
$ws4.Range("F11").Value = 7045
$ws4.Range("F12").Value = 188
$ws4.Range("F13").Value = 444
$ws4.Range("F14").Value = 9
$ws4.Range("F17").Value = 102
